$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 177 - shifts existing rows 177-186 down to 178-187
# (new weekly price entry added at the top of this variety/quality block)
$ws.Rows.Item(177).Insert()

# Populate the newly inserted row 177 with this week's data.
# Non-numeric/unchanged fields carry over the same values as the row that
# used to occupy position 177 (same market/product/category/variety/quality/unit/origin).
$ws.Range("A177").Value = 10
$ws.Range("B177").Value = "Vega Modelo de Temuco"
$ws.Range("C177").Value = "La Araucanía"
$ws.Range("D177").Value = 44585
$ws.Range("D177").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E177").Value = 9
$ws.Range("F177").Value = "Fruta"
$ws.Range("G177").Value = 100102
$ws.Range("H177").Value = "Cítricos"
$ws.Range("I177").Value = 100102006
$ws.Range("J177").Value = "Pomelo"
$ws.Range("K177").Value = "Start Ruby"
$ws.Range("L177").Value = "Primera"
$ws.Range("M177").Value = 100
$ws.Range("N177").Value = 14000
$ws.Range("O177").Value = 14000
$ws.Range("P177").Value = 14000
$ws.Range("Q177").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R177").Value = "Región de O'Higgins"
$ws.Range("S177").Value = 933
$ws.Range("T177").Value = 15
